# Update loading_percent values on Sheet1 (columns B, C, E, F, G, J, K, O; rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @"
2,10.11255850033214,5.413742889184057,22.88681357282138,41.32504043604913,3.646568668684891,8.214536716265343,9.353398781795663,22.95833993011491
3,9.790149653608982,5.239747018743675,22.52219292162736,41.18168449363715,3.648369718593737,8.241579182600864,9.121013500172868,23.08943834532558
4,9.587942994821271,5.129340549463122,22.30042578266469,41.10545222220654,3.649532512501084,8.259119109567646,8.976228469199471,23.17540166011446
5,9.50459950103104,5.083499355259263,22.21069839050067,41.07737153550131,3.650020725065171,8.26650259964069,8.916780932607587,23.21180553731609
6,9.490707295639305,5.075837654791759,22.19584146306916,41.07288953084274,3.650102661331159,8.267742881509722,8.906885210704015,23.21793325206886
7,9.58682263939156,5.12872569166054,22.29921292832136,41.10506140821683,3.649539038484632,8.259217730339744,8.975428435964291,23.17588705831905
8,10.0023452404916,5.354518282286996,22.76072235401742,41.27317784563484,3.647177880191006,8.223667062570261,9.273755677620578,23.00240652631337
9,10.77850301457844,5.767049130594259,23.67751190560336,41.69532187668209,3.642997354303535,8.161352965975215,9.83875643324893,22.70568884096945
10,11.31919620628894,6.049457564202578,24.35166184952227,42.06014804888797,3.640197104005125,8.120048797123779,10.23749230335654,22.5143186444184
11,11.55769011811444,6.173051572578441,24.65712430428887,42.2375505187474,3.638981450088456,8.102223737886927,10.41455041610923,22.43307009021909
12,11.64685108468564,6.219124713695113,24.7725123159233,42.30633105024852,3.638529434014151,8.095612013914314,10.48091678991759,22.40314108283749
13,11.62770102778921,6.209234884137846,24.74767577615875,42.29144740097571,3.638626414150862,8.097029825709868,10.4666548118912,22.40954950735052
14,11.56504893437772,6.176856801936397,24.66662365071118,42.24317728527358,3.638944095837813,8.101677019365077,10.42002440109636,22.43059099894018
15,11.52652061220999,6.156928545705953,24.61693671389341,42.21381774615957,3.639139768036216,8.104541549479867,10.39137147995747,22.44358875254672
16,11.3034529422346,6.041279817700188,24.33166530135898,42.0487809907549,3.640277717215665,8.121233075138161,10.22582885608928,22.5197455332593
17,11.16463707485955,5.96906329179533,24.15627080552252,41.95043698174329,3.640990686688647,8.13171947775739,10.12312137595488,22.56795462976353
18,11.08409478101658,5.927068855168885,24.05528079837669,41.89495180763473,3.641406247888761,8.137841780372321,10.06364251752329,22.5962295917211
19,11.05670697603968,5.912772617522236,24.02107222488095,41.87635218235709,3.641547892270723,8.139930295651579,10.04343645718428,22.60589675535945
20,11.17948720052892,5.976798418748794,24.17495380697732,41.96079442697039,3.640914223018386,8.130593787942896,10.1340970438797,22.56276612060136
21,11.58348319757722,6.186387028531759,24.69043918588342,42.25731226855317,3.638850559406904,8.10030827765592,10.43373983105594,22.42438783324329
22,11.84077622085314,6.319104058194211,25.02563552174448,42.46041966751169,3.637550344849329,8.081320509608373,10.62557633413644,22.33883595923246
23,11.70409429174723,6.248668840023837,24.84692591700264,42.35118009836575,3.63823986923884,8.091381073759758,10.52357361769818,22.3840484643404
24,11.17277574273643,5.973302848951302,24.16650769133359,41.95610853916406,3.640948774592395,8.131102420929944,10.12913628219326,22.5651101047897
25,10.57332786329006,5.658934263688852,23.42892396263294,41.57138113647213,3.644080461752022,8.177421795216361,9.688501100938993,22.78129196250262
"@

$columns = @("B", "C", "E", "F", "G", "J", "K", "O")

foreach ($line in $data -split "`n") {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $val = [double]$parts[$i + 1]
        $ws.Range("$col$row").Value = $val
    }
}
